$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header date text (stored as text, not a real date value)
$ws.Range("B1").Value = "14/03/2023"

# Update hourly values in column B (B2:B17)
$values = @(82, 165, 218, 213, 179, 148, 166, 145, 153, 162, 155, 129, 78, 48, 23, 16)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
